# Update "想去人数" (interested-count) figures in the 展览 and 全部类型
# sheets to match the newly generated site data.
#
# Each entry maps a worksheet row to its new value in column F.

$wb = $excel.ActiveWorkbook

$exhibitionUpdates = @{
    2  = 1325
    4  = 14518
    5  = 17344
    8  = 58
    12 = 52
    16 = 41
    17 = 23
    18 = 130
    20 = 1314
    25 = 7133
    27 = 33
    28 = 1159
    30 = 5835
    31 = 55
    32 = 35
    36 = 5026
}

$allTypesUpdates = @{
    2  = 1325
    4  = 14518
    5  = 17344
    8  = 58
    12 = 52
    16 = 41
    17 = 23
    18 = 130
    20 = 1314
    26 = 7133
    28 = 33
    29 = 1159
    32 = 5835
    33 = 55
    34 = 35
    38 = 5026
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
